$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: explanation text changes
$ws.Range("D3").Value = "No data"

# Row 4: status + explanation change
$ws.Range("C4").Value = "Failed"
$ws.Range("D4").Value = "No data"

# Row 6: check text + target_value in explanation
$ws.Range("A6").Value = "monthly data for March 2025 should be less than 1200"
$ws.Range("D6").Value = "tool_value=6190.0, target_value=1200.0, comparator=lt, period=month:2025-03, metric=None"

# Row 7: target_value in explanation
$ws.Range("D7").Value = "tool_value=5850.0, target_value=450.0, comparator=eq, period=month:2025-04, metric=csr_supply"

# Row 8: status + explanation (target_value and comparator) change
$ws.Range("C8").Value = "Failed"
$ws.Range("D8").Value = "tool_value=6190.0, target_value=1000.0, comparator=eq, period=month:2025-03, metric=None"

# Row 10: new row appended
$ws.Range("A10").Value = "is email address is valid"
$ws.Range("B10").Value = "evan.dummy@starlink.com"
$ws.Range("C10").Value = "Success"
$ws.Range("D10").Value = "Email found in contacts"
